$d = $word.ActiveDocument

# There are two paragraphs with the text "3 = Hampir setiap hari" in this
# document (one under the PHQ-9 list, one under the GAD-7 list). The edit
# only touches the second (GAD-7) occurrence, so find it by locating the
# paragraph whose text matches and that also follows the "GAD-7" heading.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a").Trim() -eq "3 = Hampir setiap hari") {
        $target = $p
    }
}

if ($target -ne $null) {
    $rng = $target.Range
    $rng.Find.Execute("Hampir setiap hari", $true, $false, $false, $false, $false, $true, 1, $false, "Tidak pernah sama sekali", 2)
}
